$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Swap styles of D1/E1 (content also moves) -------------------------
# Before: D1 = "Fixed" (style "Fixed"-col), E1 = "Fixed by" (style "Added by/Fixed by"-col), F1 = "Added by" (same style as E1)
# After:  D1 = "Added by" (style of old F1/E1), E1 = "Fixed" (style of old D1/C1), F1 = "Fixed by" (style of old F1/E1)
# Grab format sources before anything changes.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 5 ("centered, thin border")
$ws.Range("F1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> style 6

$ws.Range("D1").Value = "Added by"
$ws.Range("E1").Value = "Fixed"
$ws.Range("F1").Value = "Fixed by"

# --- New shrink-to-fit styles -------------------------------------------
$ws.Range("B1").ShrinkToFit = $true
$ws.Range("B2").ShrinkToFit = $true

# --- New data rows --------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Create new discussion by SUPER USER is currently not supported. Thus, all tests regarding to discussion might fail (most of then creates discussions by SUPER USER)"
$ws.Range("C3").Value = "MemberTests (maybe more)"
$ws.Range("D3").Value = "Asa"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Add moderator: in SubForum->addModerator->db.Entry fails."
$ws.Range("C4").Value = "DB (I guess)"
$ws.Range("D4").Value = "Asa"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "EditDiscussion: in Discussion->editDiscussion->db.Entry fails."
$ws.Range("C5").Value = "DB (I guess)"
$ws.Range("D5").Value = "Asa"

# --- Column widths ---------------------------------------------------------
$ws.Range("B1").EntireColumn.ColumnWidth = 89.33333333333333
$ws.Range("C1").EntireColumn.ColumnWidth = 24.333333333333332

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.75

$ws.Range("B6").Select()
